# "adding averages and more checks"
# Update the Training Dashboard's "PERIOD TO EXPIRE" (H) / "LAST UPDATE" (I)
# columns to reflect a later check-in date, tidy the Exam Dashboard's remark
# text/column width, and restyle the title + header rows (bold white text).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# --- Training Dashboard: refresh "days to expire" + "last update" ------
$newPeriods = @(400,321,324,358,360,672,400,379,377,399,357,398,401,406,405,363,133,174,177,189)

$lastUpdateRange = $ws1.Range("I3:I22")
$lastUpdateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newPeriods.Length; $i++) {
    $row = 3 + $i
    $ws1.Cells.Item($row, 8).Value = $newPeriods[$i]
    $ws1.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# --- Exam Dashboard: clearer remark + wider column for it --------------
$ws2.Range("E3").Value = "date is valid"
$ws2.Columns.Item(5).ColumnWidth = 14.17

# --- Title & header restyle: bold white text on both sheets ------------
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Bold = $true
$ws1.Range("A2:K2").Font.Color = 16777215

$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Bold = $true
$ws2.Range("A2:G2").Font.Color = 16777215
